$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.804.82"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.809.29"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.69"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.72"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.77%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.518"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.21%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.453"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.40"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000250"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "35.99"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.449.22"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.821.70"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.49"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.831.12"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.113"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.64%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.08"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "464.05"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.86"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.703"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000147"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.19"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.12"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.36%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.65%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.959.70"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.64%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.42"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.21"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.40"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.82%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.07"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0999"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.34%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.48%  "
$ws.Range("B38").Value = "Mantle"
$ws.Range("C38").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.996"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.12%  "
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.80"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.58%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.69%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "45.43"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "47.81"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.67%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.02"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.64%  "
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.40"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +12.24%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "151.27"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.96%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.36"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.86"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "394.42"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.16%  "
